# Edit script: applies the changes described by the diff to the document.
#
# Summary of changes:
#  1. The empty paragraph right after the repo hyperlink (currently centered)
#     becomes justified ("both").
#  2. The three empty, justified paragraphs that followed it are removed.
#  3. The seven body paragraphs (Dataset/Database/Pandas/Analysis/Wordcloud/
#     Preprocessing/Feature Engineering) get justified ("both") alignment.
#  4. "It is  very obvious there is..." -> "It is obvious there is..." in the
#     Analysis paragraph.
#  5. The "We also tried to add tweet length..." / "...as a feature." text in
#     the Feature Engineering paragraph is normalized into a single run
#     (the stray lastRenderedPageBreak render-cache marker goes away with it).

$d = $word.ActiveDocument

$wdAlignParagraphJustify = 3

# --- 1. Center -> Justify for the empty paragraph following the hyperlink line.
$d.Paragraphs.Item(5).Alignment = $wdAlignParagraphJustify

# --- 2. Remove the three empty justified paragraphs that come right after it.
$delStart = $d.Paragraphs.Item(6).Range.Start
$delEnd = $d.Paragraphs.Item(8).Range.End
$d.Range($delStart, $delEnd).Delete()

# --- 3. Justify the seven body paragraphs (now items 6..12 after the deletion).
for ($i = 6; $i -le 12; $i++) {
    $d.Paragraphs.Item($i).Alignment = $wdAlignParagraphJustify
}

# --- 4. Fix the "is  very obvious" -> "is obvious" wording in the Analysis paragraph.
$d.Content.Find.Execute("is  very obvious", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "is obvious", 2) | Out-Null

# --- 5. Normalize the split "We also tried ... / ... as a feature." text into
#        one contiguous run (drops the lastRenderedPageBreak cache marker too).
#        Only the back half of the text (which held the page-break marker) is
#        touched so the preceding "max_features limit of 2000." run is left
#        exactly as-is, and it merges cleanly with the "We also tried ..." run
#        that immediately precedes it.
$featPara = $d.Paragraphs.Item(12)
$featText = $featPara.Range.Text
$tailIdx = $featText.IndexOf("did not find")
$tailStart = $featPara.Range.Start + $tailIdx
$tailRange = $d.Range($tailStart, $featPara.Range.End)
$tailRange.Text = "did not find any correlation with the sentiment, so dropped the idea to use it as a feature."
